$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.294.34'
$ws.Range('E2').Value = '  -2.77%  '
$ws.Range('D3').Value = '1.572.30'
$ws.Range('E3').Value = '  -3.76%  '
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '207.86'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -3.19%  '
$ws.Range('E6').Value = '  -0.12%  '
$ws.Range('E7').Value = '  -4.75%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.244'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -2.46%  '
$ws.Range('E9').Value = '  -2.01%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '17.98'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -2.09%  '
$ws.Range('E11').Value = '  -1.14%  '
$ws.Range('D12').Value = '1.790.83'
$ws.Range('E12').Value = '  -3.75%  '
$ws.Range('D13').Value = '1.575.75'
$ws.Range('E13').Value = '  -3.64%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '4.03'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -3.32%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.505'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -3.80%  '
$ws.Range('D16').Value = '25.306.98'
$ws.Range('E16').Value = '  -2.64%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '59.73'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -2.69%  '
$ws.Range('D18').Value = '0.0₃0709'
$ws.Range('E18').Value = '  -4.47%  '
$ws.Range('E19').Value = '  -0.14%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '185.19'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -2.80%  '
$ws.Range('E21').Value = '  -2.45%  '
$ws.Range('E22').Value = '  -3.16%  '
$ws.Range('E23').Value = '  -3.26%  '
$ws.Range('E24').Value = '  -0.12%  '
$ws.Range('E25').Value = '  -2.77%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '141.05'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -1.75%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.69'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -5.48%  '
$ws.Range('E28').Value = '  -4.30%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '14.87'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -1.93%  '
$ws.Range('E30').Value = '  -6.27%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.0462'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -3.94%  '
$ws.Range('E32').Value = '  -2.60%  '
$ws.Range('E33').Value = '  -3.33%  '
$ws.Range('E34').Value = '  -1.67%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.26'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -6.16%  '
$ws.Range('D36').Value = '1.088.95'
$ws.Range('E36').Value = '  -3.84%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.00'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -0.58%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.32'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -4.90%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0150'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -2.47%  '
$ws.Range('E40').Value = '  -9.26%  '
$ws.Range('E41').Value = '  -4.47%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.756'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -2.33%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '92.99'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -5.43%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '5.07'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -2.99%  '
$ws.Range('D45').Value = '1.704.55'
$ws.Range('E45').Value = '  -3.73%  '
$ws.Range('E46').Value = '  -2.46%  '
$ws.Range('E47').Value = '  -3.54%  '
$ws.Range('E48').Value = '  -3.55%  '
$ws.Range('E49').Value = '  -4.48%  '
$ws.Range('E50').Value = '  -1.58%  '
$ws.Range('E51').Value = '  -0.22%  '
